$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.382680535316467
$ws.Range("B1").Value = 1.790967464447021
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.837021470069885
$ws.Range("E1").Value = 0.7468041181564331
